$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Price cells to keep their original text-number formatting
# (e.g. "1.00", "0.551") instead of being normalized to a floating point number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.308.34"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.171.61"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "602.02"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "154.09"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("D9").Value = "3.171.78"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").Value = "5.56"
$ws.Range("E11").Value = "  -9.98%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").Value = "38.56"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "3.691.34"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "66.374.66"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "7.39"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "3.177.54"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "510.84"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "15.37"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").Value = "0.730"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "8.11"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "14.78"
$ws.Range("E24").Value = "  -2.88%  "
$ws.Range("D25").Value = "84.61"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").Value = "9.12"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("E29").Value = "  +6.58%  "
$ws.Range("D30").Value = "3.07"
$ws.Range("E30").Value = "  +6.77%  "
$ws.Range("D31").Value = "6.90"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "27.95"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").Value = "6.51"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").Value = "509.17"
$ws.Range("D37").Value = "54.79"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "0.0887"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("E40").Value = "  +7.14%  "
$ws.Range("D41").Value = "8.82"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").Value = "0.0₃0676"
$ws.Range("E42").Value = "  +4.66%  "
$ws.Range("D43").Value = "0.297"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "2.77"
$ws.Range("E44").Value = "  -7.72%  "
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").Value = "2.835.98"
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("D47").Value = "28.03"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "2.37"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").Value = "2.58"
$ws.Range("E51").Value = "  +6.44%  "

# Restore default (unstyled) cell style now that the text values are locked in,
# so styling matches the original workbook (only the text content changed).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
